$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3916
$ws.Range("B2").Value = 2780
$ws.Range("C2").Value = 6592
$ws.Range("D2").Value = 6061
$ws.Range("E2").Value = 6957
$ws.Range("F2").Value = 5029
